$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The worksheet is protected; unprotect to allow edits, then restore protection afterward.
$ws.Unprotect()

# Update the "as of" date in the confidential disclosure note (A58).
$ws.Range("A58").Value = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution.`nModel holdings provided as of 2021-03-30 for illustrative purposes only and are subject to change."

# Refresh performance figures (columns D = trailing value, E = period return) for rows 2-55.
$ws.Range("D2").Value = 0.01628697729772464
$ws.Range("E2").Value = 0.003678658371652821
$ws.Range("D3").Value = 0.05018517985947626
$ws.Range("E3").Value = -0.006645576822412846
$ws.Range("D4").Value = 0.01496452414754223
$ws.Range("E4").Value = -0.01548291427699133
$ws.Range("D5").Value = 0.009848645545343666
$ws.Range("E5").Value = -0.01418157720344604
$ws.Range("D6").Value = 0.01606699810746234
$ws.Range("E6").Value = -0.01567571300324566
$ws.Range("D7").Value = 0.02109982567283884
$ws.Range("E7").Value = -0.01541033655593516
$ws.Range("D8").Value = 0.004354354439023813
$ws.Range("E8").Value = -0.003777148253069074
$ws.Range("D9").Value = 0.006743287327666717
$ws.Range("E9").Value = -0.01236933797909412
$ws.Range("D10").Value = 0.01400847454998778
$ws.Range("E10").Value = 0.00763150722267647
$ws.Range("D11").Value = 0.00906610571705559
$ws.Range("E11").Value = -0.004319343459794123
$ws.Range("D12").Value = 0.01455954836484603
$ws.Range("E12").Value = 0.0252824098977944
$ws.Range("D13").Value = 0.002927345091600575
$ws.Range("E13").Value = 0.03979711275848619
$ws.Range("D14").Value = 0.00607186293371229
$ws.Range("E14").Value = 0.02276086313922554
$ws.Range("D15").Value = 0.01444898769185854
$ws.Range("E15").Value = 0.01178936337437775
$ws.Range("D16").Value = 0.01051142218890208
$ws.Range("E16").Value = 0.01927912824811395
$ws.Range("D17").Value = 0.02162345511757239
$ws.Range("E17").Value = -0.003870967741935405
$ws.Range("D18").Value = 0.008781349981374115
$ws.Range("E18").Value = -0.01186202590916174
$ws.Range("D19").Value = 0.01733883744306282
$ws.Range("E19").Value = -0.006203324500120511
$ws.Range("D20").Value = 0.01216167402862385
$ws.Range("E20").Value = -0.004507888805409532
$ws.Range("D21").Value = 0.007379761756813115
$ws.Range("E21").Value = 0.005562827225130906
$ws.Range("D22").Value = 0.0139434369404057
$ws.Range("E22").Value = -0.002883355176933278
$ws.Range("D23").Value = 0.01998139416596156
$ws.Range("E23").Value = -0.008907311050357358
$ws.Range("D24").Value = 0.009863885165852862
$ws.Range("E24").Value = 0.007675333397294315
$ws.Range("D25").Value = 0.02105364994902618
$ws.Range("E25").Value = 0.008946562424438032
$ws.Range("D26").Value = 0.01152252369180774
$ws.Range("E26").Value = 0.002973712382538407
$ws.Range("D27").Value = 0.02047131370066084
$ws.Range("E27").Value = 0.003980431008784979
$ws.Range("D28").Value = 0.05545851278492295
$ws.Range("E28").Value = -0.01227448719004853
$ws.Range("D29").Value = 0.021385723564233
$ws.Range("E29").Value = 0.004608294930875667
$ws.Range("D30").Value = 0.03103661300618109
$ws.Range("E30").Value = -0.004289859948690022
$ws.Range("D31").Value = 0.01584982535695121
$ws.Range("E31").Value = -0.008836748685914553
$ws.Range("D32").Value = 0.01353271774612573
$ws.Range("E32").Value = 0.01044625675799526
$ws.Range("D33").Value = 0.02034499127883837
$ws.Range("E33").Value = -0.0005774338838203663
$ws.Range("D34").Value = 0.04005618403620004
$ws.Range("E34").Value = 0.0003323899324954027
$ws.Range("D35").Value = 0.01137978686074302
$ws.Range("E35").Value = 0.0003441156228494169
$ws.Range("D36").Value = 0.009752182337152693
$ws.Range("E36").Value = 0.0005621662138106664
$ws.Range("D37").Value = 0.01158145892658206
$ws.Range("E37").Value = 0.03338968723584101
$ws.Range("D38").Value = 0.007479488267040319
$ws.Range("E38").Value = -0.005890052356020914
$ws.Range("D39").Value = 0.01162659039374313
$ws.Range("E39").Value = 0.01774993474288689
$ws.Range("D40").Value = 0.01816719403193957
$ws.Range("E40").Value = -0.01428027418126432
$ws.Range("D41").Value = 0.01718095889054999
$ws.Range("E41").Value = 0.008011548177553118
$ws.Range("D42").Value = 0.03285198792893612
$ws.Range("E42").Value = -0.005185207259290281
$ws.Range("D43").Value = 0.01132524664115109
$ws.Range("E43").Value = -0.0008606285435229788
$ws.Range("D44").Value = 0.02170033851354598
$ws.Range("E44").Value = -0.01216726492937381
$ws.Range("D45").Value = 0.01388071627526935
$ws.Range("E45").Value = -0.03482713384960423
$ws.Range("D46").Value = 0.008200612751006538
$ws.Range("E46").Value = 0.02399541580116016
$ws.Range("D47").Value = 0.01352299310605777
$ws.Range("E47").Value = 0.01242289983494071
$ws.Range("D48").Value = 0.009707866695500682
$ws.Range("E48").Value = 0.01207115628970779
$ws.Range("D49").Value = 0.0148387564870659
$ws.Range("E49").Value = -0.01134774232707891
$ws.Range("D50").Value = 0.008577165172967175
$ws.Range("E50").Value = -0.02858806028070637
$ws.Range("D51").Value = 0.01089636129709583
$ws.Range("E51").Value = 0.01559121429384303
$ws.Range("D52").Value = 0.008871417117574078
$ws.Range("E52").Value = -0.007890264627336663
$ws.Range("D53").Value = 0.1417643670579256
$ws.Range("E53").Value = 0.0001970443349752493
$ws.Range("D54").Value = 0.04376512259849818
$ws.Range("E54").Value = -0.003422487007225139
$ws.Range("E55").Value = -0.00156281753560461

# Restore sheet protection.
$ws.Protect()
